# Generate Report for Handback
#
# Marks the handoff/handback status as complete ("Handed back: in sync with
# en-US"), fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns for both locale sheets, hyperlinks the
# newly-populated target-file cells, and widens the columns that now need to
# show the longer status text / long file names.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns + their column widths
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("E1:F1").ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# Per-locale handling (zh-cn / de-de detail sheets)
# ---------------------------------------------------------------------
$locales = @(
    @{
        Sheet = "zh-cn"
        Row2 = @{
            TargetFile = "2fcacd02-2923-41fc-9836-7569d9fdd4b6.md"
            TargetUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e7a0c8f61a1bca7cd82c45a79a2b184b5cddf89/e2e/2fcacd02-2923-41fc-9836-7569d9fdd4b6.md"
            HandbackFile = "2fcacd02-2923-41fc-9836-7569d9fdd4b6.0a73eb25752cb988a58ec15ff91181d557c6e78f.zh-cn.xlf"
            HandbackDateTime = "2016-08-16 06:25:30"
        }
        Row3 = @{
            TargetFile = "db5b626a-ada6-4089-b0fb-f74257df863b.md"
            TargetUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e7a0c8f61a1bca7cd82c45a79a2b184b5cddf89/e2e/db5b626a-ada6-4089-b0fb-f74257df863b.md"
            HandbackFile = "db5b626a-ada6-4089-b0fb-f74257df863b.d951abb2570ff2eeafbf14748d1f69c7015204b6.zh-cn.xlf"
            HandbackDateTime = "2016-08-16 06:25:30"
        }
    },
    @{
        Sheet = "de-de"
        Row2 = @{
            TargetFile = "2fcacd02-2923-41fc-9836-7569d9fdd4b6.md"
            TargetUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e7a0c8f61a1bca7cd82c45a79a2b184b5cddf89/e2e/2fcacd02-2923-41fc-9836-7569d9fdd4b6.md"
            HandbackFile = "2fcacd02-2923-41fc-9836-7569d9fdd4b6.0a73eb25752cb988a58ec15ff91181d557c6e78f.de-de.xlf"
            HandbackDateTime = "2016-08-16 06:25:37"
        }
        Row3 = @{
            TargetFile = "db5b626a-ada6-4089-b0fb-f74257df863b.md"
            TargetUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e7a0c8f61a1bca7cd82c45a79a2b184b5cddf89/e2e/db5b626a-ada6-4089-b0fb-f74257df863b.md"
            HandbackFile = "db5b626a-ada6-4089-b0fb-f74257df863b.d951abb2570ff2eeafbf14748d1f69c7015204b6.de-de.xlf"
            HandbackDateTime = "2016-08-16 06:25:37"
        }
    }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # Status column (C) for both data rows
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Widen Status (C), Latest Target File (I) and Latest Handback File (J)
    $ws.Range("C1").ColumnWidth = 29.9777047293527
    $ws.Range("I1").ColumnWidth = 40
    $ws.Range("J1").ColumnWidth = 40

    # Row 2
    $r2 = $locale.Row2
    $ws.Range("I2").Value = $r2.TargetFile
    $ws.Hyperlinks.Add($ws.Range("I2"), $r2.TargetUrl, "", "", $r2.TargetFile)
    $ws.Range("J2").Value = $r2.HandbackFile
    $ws.Range("K2").Value = $r2.HandbackDateTime

    # Row 3
    $r3 = $locale.Row3
    $ws.Range("I3").Value = $r3.TargetFile
    $ws.Hyperlinks.Add($ws.Range("I3"), $r3.TargetUrl, "", "", $r3.TargetFile)
    $ws.Range("J3").Value = $r3.HandbackFile
    $ws.Range("K3").Value = $r3.HandbackDateTime
}
